# Updates cryptos list data (prices & 1h volume %) scraped on
# Sat Jul 15 03:54:05 UTC 2023, including a row 28/29 identity swap
# (Stellar <-> EthereumClassic) to match the freshly-ranked source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.337.96"
$ws.Cells.Item(2, 5).Value = "  -3.43%  "
$ws.Cells.Item(3, 4).Value = "1.932.45"
$ws.Cells.Item(3, 5).Value = "  -3.74%  "
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).Value = "'249.10"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'0.7223"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -6.54%  "
$ws.Cells.Item(7, 4).Value = "'1.000"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.09%  "
$ws.Cells.Item(8, 4).Value = "'0.3276"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -8.78%  "
$ws.Cells.Item(9, 4).Value = "'27.27"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.66%  "
$ws.Cells.Item(10, 4).Value = "'0.06805"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -3.66%  "
$ws.Cells.Item(11, 4).Value = "'0.8037"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -4.25%  "
$ws.Cells.Item(12, 4).Value = "'0.08058"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.41%  "
$ws.Cells.Item(13, 4).Value = "1.931.23"
$ws.Cells.Item(13, 5).Value = "  -3.79%  "
$ws.Cells.Item(14, 4).Value = "'5.415"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -4.17%  "
$ws.Cells.Item(15, 4).Value = "'94.83"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -6.40%  "
$ws.Cells.Item(16, 4).Value = "'14.48"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.25%  "
$ws.Cells.Item(17, 4).Value = "30.318.33"
$ws.Cells.Item(17, 5).Value = "  -3.49%  "
$ws.Cells.Item(18, 4).Value = "'253.84"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -7.61%  "
$ws.Cells.Item(19, 4).Value = "'0.000007984"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.21%  "
$ws.Cells.Item(20, 4).Value = "'5.828"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.97%  "
$ws.Cells.Item(21, 4).Value = "2.184.08"
$ws.Cells.Item(21, 5).Value = "  -3.66%  "
$ws.Cells.Item(22, 4).Value = "'0.9998"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.03%  "
$ws.Cells.Item(23, 5).Value = "  +0.12%  "
$ws.Cells.Item(24, 4).Value = "'6.872"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -4.75%  "
$ws.Cells.Item(25, 4).Value = "'9.680"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'159.46"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.85%  "
$ws.Cells.Item(27, 4).Value = "'2.387"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.19%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).Value = "'19.06"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -5.34%  "
$ws.Cells.Item(29, 2).Value = "Stellar"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(29, 4).Value = "'0.1333"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -9.25%  "
$ws.Cells.Item(30, 4).Value = "'1.558"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.25%  "
$ws.Cells.Item(31, 4).Value = "'1.338"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.31%  "
$ws.Cells.Item(32, 4).Value = "'4.395"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -5.06%  "
$ws.Cells.Item(33, 4).Value = "'4.189"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.81%  "
$ws.Cells.Item(34, 4).Value = "'0.05062"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.27%  "
$ws.Cells.Item(36, 4).Value = "'0.7385"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -2.93%  "
$ws.Cells.Item(37, 4).Value = "'2.753"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.95%  "
$ws.Cells.Item(38, 4).Value = "'0.01969"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -2.45%  "
$ws.Cells.Item(39, 5).Value = "  -4.48%  "
$ws.Cells.Item(40, 4).Value = "'6.603"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.63%  "
$ws.Cells.Item(41, 4).Value = "'79.20"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.24%  "
$ws.Cells.Item(42, 4).Value = "'0.4461"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -5.93%  "
$ws.Cells.Item(43, 4).Value = "'1.993"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -9.00%  "
$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.05%  "
$ws.Cells.Item(45, 4).Value = "'0.8342"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.19%  "
$ws.Cells.Item(46, 4).Value = "'102.03"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.59%  "
$ws.Cells.Item(47, 4).Value = "'9.758"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.04%  "
$ws.Cells.Item(48, 4).Value = "'7.282"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -5.09%  "
$ws.Cells.Item(49, 4).Value = "'36.41"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.84%  "
$ws.Cells.Item(50, 5).Value = "  -0.62%  "
$ws.Cells.Item(51, 4).Value = "'0.4068"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -6.86%  "
